$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.306.66"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").Value = "2.620.50"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.53"
$ws.Range("E5").Value = "  +0.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.56"
$ws.Range("E6").Value = "  -1.17%  "
$ws.Range("E8").Value = "  +2.30%  "
$ws.Range("D9").Value = "2.620.78"
$ws.Range("E9").Value = "  +0.26%  "
$ws.Range("E10").Value = "  -0.49%  "
$ws.Range("E11").Value = "  +0.49%  "
$ws.Range("E12").Value = "  -0.51%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.348"
$ws.Range("E13").Value = "  -1.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.58"
$ws.Range("E14").Value = "  +0.52%  "
$ws.Range("D15").Value = "3.093.77"
$ws.Range("E15").Value = "  +0.19%  "
$ws.Range("E16").Value = "  -1.51%  "
$ws.Range("D17").Value = "67.289.25"
$ws.Range("E17").Value = "  +0.33%  "
$ws.Range("D18").Value = "2.617.49"
$ws.Range("E18").Value = "  +0.17%  "
$ws.Range("E19").Value = "  -1.29%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "362.98"
$ws.Range("E20").Value = "  +1.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.49"
$ws.Range("E21").Value = "  -3.93%  "
$ws.Range("E22").Value = "  -0.45%  "
$ws.Range("E23").Value = "  +3.36%  "
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "71.01"
$ws.Range("E25").Value = "  +5.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.08"
$ws.Range("E26").Value = "  -2.69%  "
$ws.Range("D27").Value = "2.762.43"
$ws.Range("E27").Value = "  +0.98%  "
$ws.Range("B28").Value = "Bittensor"
$ws.Range("C28").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "584.86"
$ws.Range("E28").Value = "  -4.10%  "
$ws.Range("E29").Value = "  -0.63%  "
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.14%  "
$ws.Range("E31").Value = "  -3.81%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.81"
$ws.Range("E32").Value = "  -2.11%  "
$ws.Range("E33").Value = "  -0.94%  "
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.126"
$ws.Range("E35").Value = "  -5.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.52"
$ws.Range("E36").Value = "  -1.89%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.88"
$ws.Range("E37").Value = "  -1.78%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "157.44"
$ws.Range("E38").Value = "  +1.86%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.14"
$ws.Range("E39").Value = "  -1.20%  "
$ws.Range("E40").Value = "  -0.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.26"
$ws.Range("E41").Value = "  -2.55%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.81"
$ws.Range("E42").Value = "  -1.38%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.57"
$ws.Range("E43").Value = "  +0.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.19"
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("E46").Value = "  -0.63%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "156.30"
$ws.Range("E47").Value = "  +0.64%  "
$ws.Range("D48").Value = "0.0₆0286"
$ws.Range("E48").Value = "  -2.98%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.73"
$ws.Range("E49").Value = "  -0.42%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.622"
$ws.Range("E50").Value = "  -0.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.56"
$ws.Range("E51").Value = "  -1.54%  "
